$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.428.88"
$ws.Range("E2").Value = "  +1.07%  "

$ws.Range("D3").Value = "1.638.49"
$ws.Range("E3").Value = "  +2.29%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("E5").Value = "  +0.02%  "

$ws.Range("D6").Value = "'304.92"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.75%  "

$ws.Range("D7").Value = "'0.3737"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -1.18%  "

$ws.Range("D8").Value = "'51.88"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.26%  "

$ws.Range("D9").Value = "'0.3619"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.00%  "

$ws.Range("D10").Value = "'1.255"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.94%  "

$ws.Range("D11").Value = "'0.08118"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.06%  "

$ws.Range("D13").Value = "'22.83"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.29%  "

$ws.Range("D14").Value = "'6.592"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.01%  "

$ws.Range("D15").Value = "'0.00001267"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.83%  "

$ws.Range("E16").Value = "  -1.71%  "

$ws.Range("D17").Value = "1.636.94"
$ws.Range("E17").Value = "  +2.36%  "

$ws.Range("D18").Value = "'94.20"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.24%  "

$ws.Range("D19").Value = "'0.06903"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.26%  "

$ws.Range("D20").Value = "'18.11"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.28%  "

$ws.Range("D21").Value = "'6.499"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.60%  "

$ws.Range("E22").Value = "  -0.05%  "

$ws.Range("D23").Value = "23.441.16"
$ws.Range("E23").Value = "  +1.15%  "

$ws.Range("D24").Value = "'12.72"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.85%  "

$ws.Range("D25").Value = "'2.421"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.11%  "

$ws.Range("D26").Value = "'3.035"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +1.44%  "

$ws.Range("D27").Value = "'21.19"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.15%  "

$ws.Range("D28").Value = "'151.56"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +1.08%  "

$ws.Range("D29").Value = "'5.329"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +1.59%  "

$ws.Range("D30").Value = "'135.51"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +1.32%  "

$ws.Range("E31").Value = "  -3.83%  "

$ws.Range("D32").Value = "1.819.16"

$ws.Range("D33").Value = "'6.728"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.34%  "

$ws.Range("D35").Value = "'0.02813"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +3.40%  "

$ws.Range("D36").Value = "'10.24"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.31%  "

$ws.Range("D37").Value = "'0.07243"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -3.15%  "

$ws.Range("D38").Value = "'0.2516"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.31%  "

$ws.Range("D39").Value = "'0.08778"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.22%  "

$ws.Range("D40").Value = "'6.060"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.25%  "

$ws.Range("D41").Value = "'1.372"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.97%  "

$ws.Range("E42").Value = "  -0.85%  "

$ws.Range("D43").Value = "'12.41"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.53%  "

$ws.Range("D44").Value = "'16.04"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +2.19%  "

$ws.Range("D45").Value = "'0.6491"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.71%  "

$ws.Range("D46").Value = "'2.324"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.48%  "

$ws.Range("E47").Value = "  +0.03%  "

$ws.Range("D48").Value = "'4.008"
$ws.Range("D48").ClearFormats()

$ws.Range("D49").Value = "'0.07967"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.25%  "

$ws.Range("D50").Value = "'128.08"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -3.02%  "

$ws.Range("D51").Value = "'1.200"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.29%  "
